$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column F (F2_Clientes_Compradores_Score) values for rows 2-13
$ws.Range("F2").Value = 101
$ws.Range("F3").Value = 101
$ws.Range("F4").Value = 101
$ws.Range("F5").Value = 101
$ws.Range("F6").Value = 99
$ws.Range("F7").Value = 99
$ws.Range("F8").Value = 99
$ws.Range("F9").Value = 99
$ws.Range("F10").Value = 99
$ws.Range("F11").Value = 99
$ws.Range("F12").Value = 99
$ws.Range("F13").Value = 99

# Update selection to F14
$ws.Range("F14").Select()
